# Piezo BOM.xlsx edit script
# - Row 5 (10k resistors): Quantity 2 -> 3; Positions "Rdv1, Rdv2" -> "Rdv1, Rdv2,Rdr"
# - Row 6: capacitor changed from 1.5uF to 1uF (new part, new URL, new comment); row made taller
# - Row 12: resistor changed from 1k to 200k (new part, new URL, new comment); hyperlink-style cleared
# - Selection moved to A6
# - Workbook set to non-concurrent calculation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: 10k surface mount resistors 0805 ---
$ws.Range("B5").Value = 3
$ws.Range("E5").Value = "Rdv1, Rdv2,Rdr"

# --- Row 6: capacitor swapped from 1.5uF to 1uF (per datasheet for Cf) ---
$ws.Range("A6").Value = "1uF surface mount capacitors 0805"
$ws.Range("C6").Value = "http://china.rs-online.com/web/p/ceramic-multilayer-capacitors/4515770/"
$ws.Range("D6").Value = "Must be 1uF"
$ws.Rows.Item(6).RowHeight = 30

# --- Row 12: resistor swapped from 1k to 200k (per datasheet for R14) ---
$ws.Range("A12").Value = "200k surface mount resistor 0805"
$ws.Range("C12").Value = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6791064/"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "has to be 200k"

# --- Workbook calculation: disable concurrent calculation ---
$excel.MultiThreadedCalculation.Enabled = $false

# --- Selection ---
$ws.Range("A6").Select()
